$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the date row (row 3) ---------------------------------------
# A3 currently holds "2023-06-25" (a text string, not a real date).  Simply
# assigning Range.Value = "2023-07-11" would let Excel's normal smart-entry
# logic reinterpret the ISO-looking text as a date serial, which would
# change both the stored type and the cell style. To keep it a plain text
# value (same as the source file) we build the text via a formula in a
# scratch cell and paste-special just the computed value back into A3 --
# that bypasses the "looks like a date" auto-conversion while still
# leaving a literal (non-formula) string behind.
$ws.Range("F1").Formula = "=""2023-07-11"""
$ws.Range("F1").Copy()
$ws.Range("A3").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B3").Value = 4000

# --- Update what used to be the "2023-06-26" row (row 4) ---------------
# It becomes the "Total" row; reuse the same paste-special trick for
# symmetry/safety (plain text, no style drift), then set its amount.
$ws.Range("F2").Formula = "=""Total"""
$ws.Range("F2").Copy()
$ws.Range("A4").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("B4").Value = 4000

# --- Remove the old standalone "Total" row (row 5) ----------------------
$ws.Rows(5).Delete()

# --- Clean up the scratch cells used for the paste-special trick --------
$ws.Range("F1:F2").Clear()
